$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Clear old error rows 22:25 (N/O columns) before rebuilding below
$ws.Range("N22:O25").ClearContents()

# Row 3
$ws.Range("I3").Value = "no default input, use standard word2vec init"

# Row 16
$ws.Range("I16").Value = "no default input, use standard word2vec init"

# Row 25
$ws.Range("A25").Value = "now try some default input: 12 +/- 15% variation"

# Row 26
$ws.Range("A26").Value = 12
$ws.Range("B26").Value = 0.15
$ws.Range("G26").Value = "skip-gram"
$ws.Range("I26").Value = "now try some default input: 12 +/- 15% variation"

# Row 27
$ws.Range("H27").Value = "epochs"
$ws.Range("I27").Value = "r1"
$ws.Range("J27").Value = "r2"
$ws.Range("K27").Value = "r3"
$ws.Range("L27").Value = "r4"
$ws.Range("M27").Value = "r5"
$ws.Range("N27").Value = "avg"
$ws.Range("O27").Value = "% error"

# Row 28
$ws.Range("H28").Value = 140
$ws.Range("I28").Value = 13.5287
$ws.Range("J28").Value = 13.4192
$ws.Range("K28").Value = 11.3794
$ws.Range("L28").Value = 10.9367
$ws.Range("M28").Value = 10.5267
$ws.Range("N28").Formula = "=AVERAGE(I28:M28)"
$ws.Range("O28").Formula = "=(ABS(N28-B$7)/B$7)*100"

# Row 29
$ws.Range("H29").Value = 500
$ws.Range("I29").Value = 13.3176
$ws.Range("J29").Value = 12.2884
$ws.Range("K29").Value = 11.5552
$ws.Range("L29").Value = 11.4298
$ws.Range("M29").Value = 13.3846
$ws.Range("N29").Formula = "=AVERAGE(I29:M29)"
$ws.Range("O29").Formula = "=(ABS(N29-B$7)/B$7)*100"

# Row 30
$ws.Range("H30").Value = 1000
$ws.Range("I30").Value = 13.3288
$ws.Range("J30").Value = 12.9422
$ws.Range("K30").Value = 12.5851
$ws.Range("L30").Value = 11.8316
$ws.Range("M30").Value = 13.5571
$ws.Range("N30").Formula = "=AVERAGE(I30:M30)"
$ws.Range("O30").Formula = "=(ABS(N30-B$7)/B$7)*100"

# Row 31
$ws.Range("H31").Value = 5000
$ws.Range("I31").Value = 12.1494
$ws.Range("J31").Value = 11.9102
$ws.Range("K31").Value = 12.573
$ws.Range("L31").Value = 12.3009
$ws.Range("M31").Value = 12.1928
$ws.Range("N31").Formula = "=AVERAGE(I31:M31)"
$ws.Range("O31").Formula = "=(ABS(N31-B$7)/B$7)*100"

# Row 32
$ws.Range("H32").Value = 10000
$ws.Range("I32").Value = 11.2393
$ws.Range("J32").Value = 11.7192
$ws.Range("K32").Value = 11.3957
$ws.Range("L32").Value = 11.2509
$ws.Range("M32").Value = 11.6552
$ws.Range("N32").Formula = "=AVERAGE(I32:M32)"
$ws.Range("O32").Formula = "=(ABS(N32-B$7)/B$7)*100"

# Row 33
$ws.Range("H33").Value = 25000
$ws.Range("I33").Value = 9.6038
$ws.Range("J33").Value = 9.68062
$ws.Range("K33").Value = 9.55635
$ws.Range("L33").Value = 9.58738
$ws.Range("M33").Value = 9.55759
$ws.Range("N33").Formula = "=AVERAGE(I33:M33)"
$ws.Range("O33").Formula = "=(ABS(N33-B$7)/B$7)*100"

# Row 34
$ws.Range("H34").Value = 750
$ws.Range("I34").Value = 12.2371
$ws.Range("J34").Value = 11.6394
$ws.Range("K34").Value = 13.4619
$ws.Range("L34").Value = 12.8542
$ws.Range("M34").Value = 12.6899
$ws.Range("N34").Formula = "=AVERAGE(I34:M34)"
$ws.Range("O34").Formula = "=(ABS(N34-B$7)/B$7)*100"

# Row 35
$ws.Range("H35").Value = 900
$ws.Range("I35").Value = 11.682
$ws.Range("J35").Value = 13.4823
$ws.Range("K35").Value = 13.1407
$ws.Range("L35").Value = 13.3096
$ws.Range("M35").Value = 12.896
$ws.Range("N35").Formula = "=AVERAGE(I35:M35)"
$ws.Range("O35").Formula = "=(ABS(N35-B$7)/B$7)*100"

# Row 36
$ws.Range("H36").Value = 950
$ws.Range("I36").Value = 12.0997
$ws.Range("J36").Value = 13.8833
$ws.Range("K36").Value = 13.4731
$ws.Range("L36").Value = 13.0042
$ws.Range("M36").Value = 13.0187
$ws.Range("N36").Formula = "=AVERAGE(I36:M36)"
$ws.Range("O36").Formula = "=(ABS(N36-B$7)/B$7)*100"

# Row 39
$ws.Range("G39").Value = "cbow"
$ws.Range("I39").Value = "now try some default input: 12 +/- 15% variation"

# Row 40
$ws.Range("H40").Value = "epochs"
$ws.Range("I40").Value = "r1"
$ws.Range("J40").Value = "r2"
$ws.Range("K40").Value = "r3"
$ws.Range("L40").Value = "r4"
$ws.Range("M40").Value = "r5"
$ws.Range("N40").Value = "avg"
$ws.Range("O40").Value = "% error"

# Row 41
$ws.Range("H41").Value = 4900
$ws.Range("I41").Value = 1.43904
$ws.Range("J41").Value = 1.31037
$ws.Range("K41").Value = 1.33283
$ws.Range("L41").Value = 1.4043
$ws.Range("M41").Value = 1.34423
$ws.Range("N41").Formula = "=AVERAGE(I41:M41)"
$ws.Range("O41").Formula = "=(ABS(N41-B$7)/B$7)*100"

# Row 42
$ws.Range("H42").Value = 10000
$ws.Range("I42").Value = 1.98106
$ws.Range("J42").Value = 2.00306
$ws.Range("K42").Value = 2.07386
$ws.Range("L42").Value = 1.93519
$ws.Range("M42").Value = 1.99722
$ws.Range("N42").Formula = "=AVERAGE(I42:M42)"
$ws.Range("O42").Formula = "=(ABS(N42-B$7)/B$7)*100"

# Row 43
$ws.Range("H43").Value = 25000
$ws.Range("I43").Value = 2.89062
$ws.Range("J43").Value = 2.82875
$ws.Range("K43").Value = 2.81022
$ws.Range("L43").Value = 2.85103
$ws.Range("M43").Value = 2.85029
$ws.Range("N43").Formula = "=AVERAGE(I43:M43)"
$ws.Range("O43").Formula = "=(ABS(N43-B$7)/B$7)*100"

# Row 44
$ws.Range("H44").Value = 100000
$ws.Range("I44").Value = 3.98076
$ws.Range("J44").Value = 3.95892
$ws.Range("K44").Value = 3.98268
$ws.Range("L44").Value = 3.95722
$ws.Range("M44").Value = 3.94627
$ws.Range("N44").Formula = "=AVERAGE(I44:M44)"
$ws.Range("O44").Formula = "=(ABS(N44-B$7)/B$7)*100"

# Row 45
$ws.Range("H45").Value = 250000
$ws.Range("I45").Value = 4.60704
$ws.Range("J45").Value = 4.60674
$ws.Range("K45").Value = 4.59767
$ws.Range("L45").Value = 4.58531
$ws.Range("M45").Value = 4.59661
$ws.Range("N45").Formula = "=AVERAGE(I45:M45)"
$ws.Range("O45").Formula = "=(ABS(N45-B$7)/B$7)*100"

# Row 48
$ws.Range("A48").Value = "try default sticky = 10 (instead of 1), no default init"

# Row 49
$ws.Range("G49").Value = "skip-gram"
$ws.Range("I49").Value = "default sticky = 10"

# Row 50
$ws.Range("H50").Value = "epochs"
$ws.Range("I50").Value = "r1"
$ws.Range("J50").Value = "r2"
$ws.Range("K50").Value = "r3"
$ws.Range("L50").Value = "r4"
$ws.Range("M50").Value = "r5"
$ws.Range("N50").Value = "avg"
$ws.Range("O50").Value = "% error"

# Row 51
$ws.Range("H51").Value = 140
$ws.Range("I51").Value = 3.94842
$ws.Range("J51").Value = 4.20184
$ws.Range("K51").Value = 4.14311
$ws.Range("L51").Value = 4.02965
$ws.Range("M51").Value = 4.12684
$ws.Range("N51").Formula = "=AVERAGE(I51:M51)"
$ws.Range("O51").Formula = "=(ABS(N51-B$7)/B$7)*100"

# Row 52
$ws.Range("H52").Value = 4000
$ws.Range("I52").Value = 10.9821
$ws.Range("J52").Value = 11.0712
$ws.Range("K52").Value = 10.8311
$ws.Range("L52").Value = 11.0284
$ws.Range("M52").Value = 10.8658
$ws.Range("N52").Formula = "=AVERAGE(I52:M52)"
$ws.Range("O52").Formula = "=(ABS(N52-B$7)/B$7)*100"

# Row 53
$ws.Range("H53").Value = 10000
$ws.Range("I53").Value = 10.768
$ws.Range("J53").Value = 10.5092
$ws.Range("K53").Value = 10.6591
$ws.Range("L53").Value = 10.6387
$ws.Range("M53").Value = 10.5443
$ws.Range("N53").Formula = "=AVERAGE(I53:M53)"
$ws.Range("O53").Formula = "=(ABS(N53-B$7)/B$7)*100"

# Row 54
$ws.Range("H54").Value = 6000
$ws.Range("I54").Value = 10.9183
$ws.Range("J54").Value = 11.0059
$ws.Range("K54").Value = 10.9774
$ws.Range("L54").Value = 10.8982
$ws.Range("M54").Value = 11.0359
$ws.Range("N54").Formula = "=AVERAGE(I54:M54)"
$ws.Range("O54").Formula = "=(ABS(N54-B$7)/B$7)*100"

# Row 55
$ws.Range("H55").Value = 7000
$ws.Range("I55").Value = 10.7441
$ws.Range("J55").Value = 10.954
$ws.Range("K55").Value = 10.9614
$ws.Range("L55").Value = 10.7574
$ws.Range("M55").Value = 11.0132
$ws.Range("N55").Formula = "=AVERAGE(I55:M55)"
$ws.Range("O55").Formula = "=(ABS(N55-B$7)/B$7)*100"

# Row 57
$ws.Range("G57").Value = "cbow"
$ws.Range("I57").Value = "default sticky = 10"

# Row 58
$ws.Range("H58").Value = "epochs"
$ws.Range("I58").Value = "r1"
$ws.Range("J58").Value = "r2"
$ws.Range("K58").Value = "r3"
$ws.Range("L58").Value = "r4"
$ws.Range("M58").Value = "r5"
$ws.Range("N58").Value = "avg"
$ws.Range("O58").Value = "% error"

# Row 59
$ws.Range("H59").Value = 4900
$ws.Range("I59").Value = 1.6166
$ws.Range("J59").Value = 1.411
$ws.Range("K59").Value = 1.4975
$ws.Range("L59").Value = 1.70159
$ws.Range("M59").Value = 1.40573
$ws.Range("N59").Formula = "=AVERAGE(I59:M59)"
$ws.Range("O59").Formula = "=(ABS(N59-B$7)/B$7)*100"

# Row 60
$ws.Range("H60").Value = 20000
$ws.Range("I60").Value = 4.15007
$ws.Range("J60").Value = 4.01731
$ws.Range("K60").Value = 4.07246
$ws.Range("L60").Value = 4.07214
$ws.Range("M60").Value = 3.93264
$ws.Range("N60").Formula = "=AVERAGE(I60:M60)"
$ws.Range("O60").Formula = "=(ABS(N60-B$7)/B$7)*100"

# Row 61
$ws.Range("H61").Value = 100000
$ws.Range("I61").Value = 5.28359
$ws.Range("J61").Value = 5.27153
$ws.Range("K61").Value = 5.28345
$ws.Range("L61").Value = 5.27464
$ws.Range("M61").Value = 5.27037
$ws.Range("N61").Formula = "=AVERAGE(I61:M61)"
$ws.Range("O61").Formula = "=(ABS(N61-B$7)/B$7)*100"

# Update view / selection to match target state
$ws.Application.ActiveWindow.ScrollRow = 47
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("J63").Select()
